# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "codcom" column changes from a curated dimension to a plain measure
$ws.Range("D2").Value = "iaest-measure:codcom"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"

# "direccion-provincial-nombre" column changes from a dimension (refArea/URI-Provincia)
# to a plain measure (xsd:int)
$ws.Range("H2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"

# Remove the now-obsolete codcom mapping file row
$ws.Rows.Item(5).Delete()
